# Apply the edit described by the diff:
# - Delete the first 8 data rows (old rows 2-9), shifting rows 10-22 up to rows 2-14
# - Append 7 brand-new data rows after the shifted data (new rows 15-21)
# - The sheet dimension shrinks from A1:H22 to A1:H21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old data rows 2 through 9 (inclusive); this shifts rows 10-22 up to 2-14
$ws.Range("A2:H9").EntireRow.Delete() | Out-Null

# New rows to append after the shifted data (will land at rows 15-21)
$newRows = @(
    @(1300, 'struggle', -11.7663733065128, -6.984511554241188, 3.816117525100708, -1.965306162834168, -0.784503698348999, 1.217912554740906),
    @(1400, 'struggle', -1.357394456863403, -10.09910678863525, 3.819830894470215, -0.3843869566917419, -3.860210180282593, 2.725528001785278),
    @(1500, 'struggle', 1.284981921315195, -13.9884957075119, -13.69542229175569, 4.380514621734619, -3.570049285888672, 1.001513600349426),
    @(1600, 'struggle', 1.272318005561828, -9.928469419479365, -8.66020488739013, 0.96409809589386, -0.7756461501121521, -0.2683225572109222),
    @(1700, 'struggle', 2.418770149350169, -4.946553826332086, -8.305895447731023, -0.6479753255844116, 0.6890559792518616, 1.091005325317383),
    @(1800, 'struggle', -0.9577411413192878, -7.853628158569351, -5.728095054626454, -1.356426239013672, 3.433979034423828, -1.384373307228088),
    @(1900, 'struggle', -5.118649840354919, -10.2695299386978, -1.659017741680144, 0.2755002379417419, 2.776687860488892, -1.657124638557434)
)

$startRow = 15
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
